$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 100

$ws.Range("C12").Value = 0.2632
$ws.Range("D12").Value = 0.2105
$ws.Range("E12").Value = 0.5263

$ws.Range("C13").Value = 0.303
$ws.Range("D13").Value = 0.0909
$ws.Range("E13").Value = 0.6061

$ws.Range("C14").Value = 0.3191
$ws.Range("D14").Value = 0.0426
$ws.Range("E14").Value = 0.6383

$ws.Range("Q23").Formula = "=E12+D12"
$ws.Range("D26").Formula = "=E13+D13"
